$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values computed to replace the old Strike# based values,
# for rows 2-17 (column G). Row 18 is unchanged (already 0).
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 3
    9  = 0
    10 = 3
    11 = 2
    12 = 0
    13 = 2
    14 = 2
    15 = 2
    16 = 0
    17 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
